$d = $word.ActiveDocument

function Split-ReplaceText {
    param(
        [string]$oldText,
        [string[]]$segments
    )
    $searchRange = $d.Content
    $found = $searchRange.Find.Execute($oldText, $true, $false, $false, $false, $false, $true, 1, $false, "", 0)
    if (-not $found) {
        Write-Host "NOT FOUND:" $oldText
        return
    }

    $startPos = $searchRange.Start
    $newText = [string]::Join("", $segments)

    # Replace the found range's text with the concatenated replacement text.
    $searchRange.Text = $newText

    # Now split the replaced text back into separate runs matching the
    # segment boundaries, by toggling a formatting property on/off for the
    # trailing portion at each boundary (forces the engine to create a new
    # run there, instead of silently re-merging runs with identical rPr).
    $cursor = $startPos
    for ($i = 0; $i -lt $segments.Length; $i++) {
        $segLen = $segments[$i].Length
        $segEnd = $cursor + $segLen
        if ($i -lt ($segments.Length - 1)) {
            $restRange = $d.Range($segEnd, $startPos + $newText.Length)
            $restRange.Font.Bold = 1
            $restRange.Font.Bold = 0
        }
        $cursor = $segEnd
    }
}

# --- Header row abbreviations ---
Split-ReplaceText "total dissolved N" @("TD", "N")
Split-ReplaceText "dissolved reactive P (mg/l)" @("DR", "P (mg/l)")
Split-ReplaceText "total C" @("TO", "C")
Split-ReplaceText "total N" @("T", "N")
Split-ReplaceText "total P" @("T", "P")
Split-ReplaceText "total S (mg/l)" @("T", "S (mg/l)")

# --- Numeric ranges: remove the surrounding spaces around the en dash ---
Split-ReplaceText "0 – 3.15" @("0–", "3.15")
Split-ReplaceText "0 – 10.5" @("0–", "10.5")
Split-ReplaceText "0 – 2.07" @("0–", "2.07")
Split-ReplaceText "0 – 5.11" @("0", "–5.11")
Split-ReplaceText "0 – 1.55" @("0–", "1.55")
Split-ReplaceText "0 – 0.032" @("0–", "0.032")
